$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 265; this shifts existing rows 265-320
# down to 266-321 (row 321 becomes a copy of the former row 320, and
# the new empty row 265 inherits formatting, e.g. the date style on D).
$ws.Rows("265:265").Insert()

# Populate the newly inserted row 265 with the new weekly record.
$ws.Range("A265").Value = 11
$ws.Range("B265").Value = "Vega Monumental Concepción"
$ws.Range("C265").Value = "Bíobío"
$ws.Range("D265").Value = 44637
$ws.Range("E265").Value = 8
$ws.Range("F265").Value = "Fruta"
$ws.Range("G265").Value = 100101
$ws.Range("H265").Value = "Berries"
$ws.Range("I265").Value = 100112025
$ws.Range("J265").Value = "Frutilla"
$ws.Range("K265").Value = "Sin especificar"
$ws.Range("L265").Value = "Primera"
$ws.Range("M265").Value = 220
$ws.Range("N265").Value = 7000
$ws.Range("O265").Value = 8500
$ws.Range("P265").Value = 7682
$ws.Range("Q265").Value = "$/bandeja 7 kilos"
$ws.Range("R265").Value = "Región del Maule"
$ws.Range("S265").Value = 1097
$ws.Range("T265").Value = 7
